# Add a new lookup row (OC215 / QDF entry message) to the bottom of the
# problem-description table on Sheet1, mirroring the style of the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 36

$codeCell = $ws.Cells.Item($newRow, 1)   # column A
$descCell = $ws.Cells.Item($newRow, 2)   # column B

$codeCell.Value = "OC215"
$descCell.Value = "QDF entry '' must have exactly 4 characters"

# Match the formatting used by the rest of the table (column A centered
# horizontally + vertically, column B centered horizontally).
$codeCell.HorizontalAlignment = -4108 # xlCenter
$codeCell.VerticalAlignment = -4108   # xlCenter

$descCell.HorizontalAlignment = -4108 # xlCenter
